# Update dictionary counts for years 2008-2010 (rows 2-4) per the new
# data source ("new dic, counted till 2010").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - year 2008
$ws.Range("C2").Value = 172
$ws.Range("D2").Value = 611
$ws.Range("E2").Value = 111
$ws.Range("G2").Value = 13
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 6
$ws.Range("K2").Value = 60
$ws.Range("M2").Value = 65

# Row 3 - year 2009
$ws.Range("B3").Value = 1794
$ws.Range("C3").Value = 372
$ws.Range("D3").Value = 930
$ws.Range("E3").Value = 257
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 7
$ws.Range("K3").Value = 98
$ws.Range("M3").Value = 54

# Row 4 - year 2010
$ws.Range("C4").Value = 412
$ws.Range("D4").Value = 1067
$ws.Range("E4").Value = 467
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 5
$ws.Range("K4").Value = 87
$ws.Range("M4").Value = 83
